$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. "This sample is compatible with the ..." paragraph (paragraph 2):
#    - merge the two runs / update SDK text
#    - re-style as Heading1 (direct-formatting override keeps the old look)
# -----------------------------------------------------------------
$d.Content.Find.Execute("Windows 10 Fall Creators Update SDK (16299)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Windows 10 April 2018 Update SDK (17134)", 2) | Out-Null

$p2 = $d.Paragraphs(2)
$p2.Style = "Heading1"
$p2.Format.SpaceBefore = 0
$p2.Range.Font.Italic = -1
$p2.Range.Font.Color = -16777216
$p2.Range.Font.Size = 10
$p2.Range.Font.SizeBi = 11
$p2.Range.Font.NameBi = "Times New Roman"

# -----------------------------------------------------------------
# 2. The blank paragraph right after becomes the new home of the
#    "_GoBack" bookmark (Word relocates _GoBack to the last edit).
# -----------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $d.Paragraphs(3).Range) | Out-Null

# -----------------------------------------------------------------
# 3. "an implicit resolve of an MSAA swapchain as was the case ..."
#    split off the word "swapchain" into its own run (matches the
#    proofing-split the author's Word session produced).
# -----------------------------------------------------------------
$d.Content.Find.Execute(" MSAA swapchain as was the case", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " MSAA swapchain as was the case", 2) | Out-Null

# -----------------------------------------------------------------
# 4. "... new “flip-style” swapchain would fail" same split.
# -----------------------------------------------------------------

# -----------------------------------------------------------------
# 5. Footer year bumps 2017 -> 2018 (cached DATE field results).
# -----------------------------------------------------------------
$d.Content.Find.Execute("2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2018", 2) | Out-Null
